$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 2777.8572
$ws.Range("I99").Value = 2289
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 6867
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -5369
$ws.Range("N99").Value = -14996

# Row 113
$ws.Range("H113").Value = 2273.8235
$ws.Range("J113").Value = 2256.6667
$ws.Range("L113").Value = 2256.6667
$ws.Range("N113").Value = -8764.6667

# Row 116
$ws.Range("H116").Value = 3727.652
$ws.Range("I116").Value = 2583.5
$ws.Range("J116").Value = 6342.857
$ws.Range("K116").Value = 2583.5
$ws.Range("L116").Value = 6342.857
$ws.Range("M116").Value = 858.5
$ws.Range("N116").Value = -13226.857

# Row 121
$ws.Range("H121").Value = 1359.8
$ws.Range("J121").Value = 2952
$ws.Range("L121").Value = 8856
$ws.Range("N121").Value = -12350

# Row 132
$ws.Range("H132").Value = 21233.375
$ws.Range("I132").Value = 3608.3235
$ws.Range("K132").Value = 10824.9705
$ws.Range("M132").Value = -8294.970499999999

# Row 137
$ws.Range("H137").Value = 3089.1292
$ws.Range("I137").Value = 965.1087
$ws.Range("J137").Value = 9195.6875
$ws.Range("K137").Value = 2895.3261
$ws.Range("L137").Value = 27587.0625
$ws.Range("M137").Value = -345.3261000000002
$ws.Range("N137").Value = -32687.0625

# Row 138
$ws.Range("H138").Value = 1660.7012
$ws.Range("J138").Value = 2350.524
$ws.Range("L138").Value = 7051.572
$ws.Range("N138").Value = -17331.572

# Row 141
$ws.Range("H141").Value = 2441.5908
$ws.Range("I141").Value = 1096.5625
$ws.Range("K141").Value = 3289.6875
$ws.Range("M141").Value = 1890.3125

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1232
$ws.Range("I61").Value = 949.9
$ws.Range("J61").Value = 1883
$ws.Range("K61").Value = 949.9
$ws.Range("L61").Value = 1883
$ws.Range("M61").Value = -737.9
$ws.Range("N61").Value = -2307

# Row 97
$ws.Range("H97").Value = 1432
$ws.Range("I97").Value = 1502.3529
$ws.Range("J97").Value = 1033.3334
$ws.Range("K97").Value = 1502.3529
$ws.Range("L97").Value = 1033.3334
$ws.Range("M97").Value = -1006.3529
$ws.Range("N97").Value = -2025.3334

# Row 122
$ws.Range("H122").Value = 1045.7
$ws.Range("I122").Value = 939.6667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2819.0001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -369.0001000000002
$ws.Range("N122").Value = -10900

# Row 132
$ws.Range("H132").Value = 3098.4443
$ws.Range("I132").Value = 1923.5555
$ws.Range("J132").Value = 4273.3335
$ws.Range("K132").Value = 5770.666499999999
$ws.Range("L132").Value = 12820.0005
$ws.Range("M132").Value = -3240.666499999999
$ws.Range("N132").Value = -17880.0005

# Row 136
$ws.Range("H136").Value = 1232
$ws.Range("I136").Value = 949.9
$ws.Range("J136").Value = 1883
$ws.Range("K136").Value = 2849.7
$ws.Range("L136").Value = 5649
$ws.Range("M136").Value = -299.6999999999998
$ws.Range("N136").Value = -10749

$ws = $wb.Worksheets.Item("BSM")
# Row 63
$ws.Range("H63").Value = 29181.428
$ws.Range("J63").Value = 29181.428
$ws.Range("L63").Value = 29181.428
$ws.Range("N63").Value = -30553.428

# Row 66
$ws.Range("H66").Value = 29181.428
$ws.Range("J66").Value = 29181.428
$ws.Range("L66").Value = 87544.284
$ws.Range("N66").Value = -94408.284

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 94
$ws.Range("H94").Value = 1444
$ws.Range("I94").Value = 1433.4
$ws.Range("J94").Value = 1550
$ws.Range("K94").Value = 1433.4
$ws.Range("L94").Value = 1550
$ws.Range("M94").Value = -982.4000000000001
$ws.Range("N94").Value = -2452

# Row 125
$ws.Range("H125").Value = 50776
$ws.Range("J125").Value = 50776
$ws.Range("L125").Value = 50776
$ws.Range("N125").Value = -60616

# Row 134
$ws.Range("H134").Value = 3372.8364
$ws.Range("I134").Value = 1907.2
$ws.Range("J134").Value = 3698.5334
$ws.Range("K134").Value = 5721.6
$ws.Range("L134").Value = 11095.6002
$ws.Range("M134").Value = -3186.6
$ws.Range("N134").Value = -16165.6002

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2971.14
$ws.Range("I31").Value = 849.6799999999999
$ws.Range("J31").Value = 3678.2932
$ws.Range("K31").Value = 849.6799999999999
$ws.Range("L31").Value = 3678.2932
$ws.Range("M31").Value = -554.6799999999999
$ws.Range("N31").Value = -4268.2932

# Row 34
$ws.Range("H34").Value = 2971.14
$ws.Range("I34").Value = 849.6799999999999
$ws.Range("J34").Value = 3678.2932
$ws.Range("K34").Value = 849.6799999999999
$ws.Range("L34").Value = 3678.2932
$ws.Range("M34").Value = -647.6799999999999
$ws.Range("N34").Value = -4082.2932

# Row 141
$ws.Range("H141").Value = 5354.4443
$ws.Range("J141").Value = 5354.4443
$ws.Range("L141").Value = 5354.4443
$ws.Range("N141").Value = -15714.4443

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3246.465
$ws.Range("I5").Value = 6171.4443
$ws.Range("J5").Value = 1140.48
$ws.Range("K5").Value = 18514.3329
$ws.Range("L5").Value = 3421.44
$ws.Range("M5").Value = -18402.3329
$ws.Range("N5").Value = -3645.44

# Row 75
$ws.Range("H75").Value = 2730
$ws.Range("I75").Value = 192.5
$ws.Range("J75").Value = 3455
$ws.Range("K75").Value = 577.5
$ws.Range("L75").Value = 10365
$ws.Range("M75").Value = 420.5
$ws.Range("N75").Value = -12361

# Row 78
$ws.Range("H78").Value = 2730
$ws.Range("I78").Value = 192.5
$ws.Range("J78").Value = 3455
$ws.Range("K78").Value = 1732.5
$ws.Range("L78").Value = 31095
$ws.Range("M78").Value = 3259.5
$ws.Range("N78").Value = -41079

# Row 92
$ws.Range("H92").Value = 1256.4445
$ws.Range("I92").Value = 950.5
$ws.Range("J92").Value = 1501.2
$ws.Range("K92").Value = 2851.5
$ws.Range("L92").Value = 4503.6
$ws.Range("M92").Value = -1603.5
$ws.Range("N92").Value = -6999.6

# Row 102
$ws.Range("H102").Value = 19900
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 19900
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 59700
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -64568

# Row 113
$ws.Range("H113").Value = 2593.4614
$ws.Range("I113").Value = 3584.7576
$ws.Range("J113").Value = 871.7368
$ws.Range("K113").Value = 10754.2728
$ws.Range("L113").Value = 2615.2104
$ws.Range("M113").Value = -8584.272799999999
$ws.Range("N113").Value = -6955.2104

# Row 123
$ws.Range("H123").Value = 2616.6667
$ws.Range("J123").Value = 2933.3333
$ws.Range("L123").Value = 8799.999899999999
$ws.Range("N123").Value = -13699.9999

# Row 131
$ws.Range("H131").Value = 3530.3333
$ws.Range("I131").Value = 7138.4
$ws.Range("J131").Value = 1525.8518
$ws.Range("K131").Value = 21415.2
$ws.Range("L131").Value = 4577.555399999999
$ws.Range("M131").Value = -16375.2
$ws.Range("N131").Value = -14657.5554

# Row 132
$ws.Range("H132").Value = 1861.5807
$ws.Range("I132").Value = 1407
$ws.Range("J132").Value = 2346.4666
$ws.Range("K132").Value = 12663
$ws.Range("L132").Value = 21118.1994
$ws.Range("M132").Value = -10133
$ws.Range("N132").Value = -26178.1994

# Row 135
$ws.Range("H135").Value = 3246.465
$ws.Range("I135").Value = 6171.4443
$ws.Range("J135").Value = 1140.48
$ws.Range("K135").Value = 55542.9987
$ws.Range("L135").Value = 10264.32
$ws.Range("M135").Value = -53007.9987
$ws.Range("N135").Value = -15334.32

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 4366.6665
$ws.Range("J21").Value = 4366.6665
$ws.Range("L21").Value = 4366.6665
$ws.Range("N21").Value = -4712.6665

# Row 30
$ws.Range("H30").Value = 4366.6665
$ws.Range("J30").Value = 4366.6665
$ws.Range("L30").Value = 4366.6665
$ws.Range("N30").Value = -4576.6665

# Row 102
$ws.Range("H102").Value = 1599.6666
$ws.Range("I102").Value = 1599.6666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1599.6666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 22.33339999999998
$ws.Range("N102").ClearContents()

# Row 126
$ws.Range("H126").Value = 3312.8667
$ws.Range("I126").Value = 2891.1667
$ws.Range("J126").Value = 4999.6665
$ws.Range("K126").Value = 8673.500100000001
$ws.Range("L126").Value = 14998.9995
$ws.Range("M126").Value = -6203.500100000001
$ws.Range("N126").Value = -19938.9995

# Row 132
$ws.Range("H132").Value = 2684.976
$ws.Range("I132").Value = 1929.75
$ws.Range("J132").Value = 4195.4287
$ws.Range("K132").Value = 5789.25
$ws.Range("L132").Value = 12586.2861
$ws.Range("M132").Value = -3259.25
$ws.Range("N132").Value = -17646.2861

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3057.9375
$ws.Range("I7").Value = 2210.6365
$ws.Range("J7").Value = 4922
$ws.Range("K7").Value = 2210.6365
$ws.Range("L7").Value = 4922
$ws.Range("M7").Value = -2098.6365
$ws.Range("N7").Value = -5146

# Row 22
$ws.Range("H22").Value = 497.7143
$ws.Range("I22").Value = 161.33333
$ws.Range("K22").Value = 161.33333
$ws.Range("M22").Value = 133.66667

# Row 27
$ws.Range("H27").Value = 497.7143
$ws.Range("I27").Value = 161.33333
$ws.Range("K27").Value = 161.33333
$ws.Range("M27").Value = -54.33332999999999

# Row 40
$ws.Range("H40").Value = 3398.4443
$ws.Range("I40").Value = 2643
$ws.Range("J40").Value = 6042.5
$ws.Range("K40").Value = 2643
$ws.Range("L40").Value = 6042.5
$ws.Range("M40").Value = -2507
$ws.Range("N40").Value = -6314.5

# Row 125
$ws.Range("H125").Value = 49707.332
$ws.Range("J125").Value = 49707.332
$ws.Range("L125").Value = 49707.332
$ws.Range("N125").Value = -59547.332

# Row 126
$ws.Range("H126").Value = 3057.9375
$ws.Range("I126").Value = 2210.6365
$ws.Range("J126").Value = 4922
$ws.Range("K126").Value = 6631.9095
$ws.Range("L126").Value = 14766
$ws.Range("M126").Value = -4161.9095
$ws.Range("N126").Value = -19706

# Row 132
$ws.Range("H132").Value = 2680.9622
$ws.Range("I132").Value = 2037
$ws.Range("J132").Value = 3743.5
$ws.Range("K132").Value = 6111
$ws.Range("L132").Value = 11230.5
$ws.Range("M132").Value = -3581
$ws.Range("N132").Value = -16290.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2942365.2
$ws.Range("I126").Value = 3677694.8
$ws.Range("J126").Value = 1047.5
$ws.Range("K126").Value = 11033084.4
$ws.Range("L126").Value = 3142.5
$ws.Range("M126").Value = -11030614.4
$ws.Range("N126").Value = -8082.5

# Row 132
$ws.Range("H132").Value = 1485.6735
$ws.Range("I132").Value = 1282.1177
$ws.Range("J132").Value = 1947.0667
$ws.Range("K132").Value = 3846.3531
$ws.Range("L132").Value = 5841.2001
$ws.Range("M132").Value = -1316.3531
$ws.Range("N132").Value = -10901.2001
